# Update latest output (run 167)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 315.7372140000001
$schedule.Range("F2").Value = 6.96069695767196
$schedule.Range("E3").Value = 440.1257055
$schedule.Range("F3").Value = 29.10884295634921

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B21").Value = -7.84163

$detailed.Range("B22").Value = -7.67981

$detailed.Range("B23").Value = -9.029170000000001
$detailed.Range("C23").Value = "historical"

$detailed.Range("B24").Value = -10
$detailed.Range("C24").Value = "historical"

$detailed.Range("B25").Value = -14
$detailed.Range("C25").Value = "historical"

$detailed.Range("B26").Value = -14
$detailed.Range("C26").Value = "historical"

$detailed.Range("B27").Value = -7.49419

$detailed.Range("B28").Value = -7.799

$detailed.Range("B29").Value = -8.495839999999999

$detailed.Range("B30").Value = -6.97469

$detailed.Range("B31").Value = -5.53787

$detailed.Range("B32").Value = -5.50985

$detailed.Range("B33").Value = -2.69285

$detailed.Range("B34").Value = -0.02892

$detailed.Range("B35").Value = -6.72418

$detailed.Range("B36").Value = -7.88262

$detailed.Range("B37").Value = -5.26657

$detailed.Range("B38").Value = 0.00015

$detailed.Range("B39").Value = 10.41456

$detailed.Range("B40").Value = 22.55224

$detailed.Range("B41").Value = 55.33036

$detailed.Range("B42").Value = 53.90468

$detailed.Range("B44").Value = 57.04367

$detailed.Range("B45").Value = 57.03541

$detailed.Range("B46").Value = 57.02642

$detailed.Range("B47").Value = 57.03043

$detailed.Range("B48").Value = 57.06

$detailed.Range("B49").Value = 56.98
